$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("E1").Value = "TrackID"
$ws.Range("J1").Value = "Acao"

# --- Row 2 (new test data) ---
$ws.Range("A2").Value = "Caio"
$ws.Range("C2").Value = "FE09"
$ws.Range("D2").Value = "FE"
$ws.Range("E2").Value = "JDKLSJ78JLK"
$ws.Range("G2").Value = "Test"
$ws.Range("H2").Value = "Test"
$ws.Range("I2").Value = "Test"
$ws.Range("J2").Value = "Test"
$ws.Range("K2").Value = "25/03/2025 11:22"

# --- Row 3 (new test data) ---
$ws.Range("A3").Value = "Benda"
$ws.Range("C3").Value = "FE09"
$ws.Range("D3").Value = "BE"
$ws.Range("E3").Value = "JKDLS89JKL"
$ws.Range("F3").Value = "Manila"
$ws.Range("G3").Value = "TESTE"
$ws.Range("H3").Value = "TESTE"
$ws.Range("I3").Value = "TESTE"
$ws.Range("J3").Value = "TESTE"
$ws.Range("K3").Value = "25/03/2025 11:24"

# --- Remove old row 4 (data shrank from 4 rows to 3) ---
$ws.Rows.Item(4).Delete()

# --- Selection moves to G17 ---
$ws.Range("G17").Select() | Out-Null
